# ------------------------------------------------------------------
# Commit: "[TEST SCRAPE] updated files from azure vm"
#
# 1) Remove the stray empty INNING_NUMBER cells (B9, B10, B16) on the
#    "ODI Batting" sheet - they were placeholders with no real value.
# 2) Add a brand-new "ODI Batting Extra" sheet (after "ODI Bowling")
#    with per-match batting-position / boundary-count / match-share
#    stats for player 5926.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) Clear out the empty placeholder cells on "ODI Batting" -----
$battingWs = $wb.Worksheets.Item("ODI Batting")
foreach ($r in @(9, 10, 16)) {
    $battingWs.Range("B$r").Value = ""
}

# --- 2) Create the new "ODI Batting Extra" sheet --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

# Data rows (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH)
$data = @(
    @{A="4437"; B=5; C="1"; D="1"; E="7.27%"; F="NO"},
    @{A="4564"; B=7; C="3"; D="1"; E="12.78%"; F="NO"},
    @{A="4565"; B=$null; C=$null; D=$null; E=$null; F="NO"},
    @{A="4567"; B=7; C="1"; D="1"; E="16.19%"; F="NO"},
    @{A="4600"; B=$null; C=$null; D=$null; E=$null; F="NO"},
    @{A="4601"; B=8; C="0"; D="0"; E="5.12%"; F="NO"},
    @{A="4603"; B=8; C="2"; D="1"; E="15.24%"; F="NO"},
    @{A="4644"; B=8; C=$null; D=$null; E=$null; F="YES"},
    @{A="4645"; B=6; C=$null; D=$null; E=$null; F="NO"},
    @{A="4646"; B=6; C="0"; D="0"; E="2.13%"; F="NO"},
    @{A="4647"; B=$null; C=$null; D=$null; E=$null; F="NO"},
    @{A="4649"; B=7; C="2"; D="2"; E="9.36%"; F="NO"},
    @{A="4660"; B=$null; C=$null; D=$null; E=$null; F="NO"},
    @{A="4725"; B=6; C="0"; D="0"; E="6.38%"; F="NO"},
    @{A="4728"; B=5; C=$null; D=$null; E=$null; F="NO"}
)

$r = 2
foreach ($row in $data) {
    # MATCH_CODE - always text, even though it looks numeric
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row.A

    # BATTING_POSITION - real number when present, otherwise a blank text cell
    if ($row.B -eq $null) {
        $ws.Range("B$r").NumberFormat = "@"
        $ws.Range("B$r").Value = ""
    } else {
        $ws.Range("B$r").Value = $row.B
    }

    # NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL - text columns, blank when absent
    foreach ($col in @("C", "D", "E")) {
        $ws.Range("$col$r").NumberFormat = "@"
        $val = $row[$col]
        if ($val -eq $null) {
            $ws.Range("$col$r").Value = ""
        } else {
            $ws.Range("$col$r").Value = $val
        }
    }

    # MAN_OF_MATCH - always present text
    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = $row.F

    $r++
}
